$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Matricula"
$ws.Range("B1").Value = "Valor"

# Update data rows (A4 reuses the same text as A2, per the diff's shared-string reuse)
$ws.Range("A2").Value = "111-x"
$ws.Range("A3").Value = "12-x"
$ws.Range("A4").Value = "111-x"

# Apply the new formatting style to A2:A4
$ws.Range("A2:A4").Font.Bold = $true
